$d = $word.ActiveDocument

# --------------------------------------------------------------------
# Helper: insert $text at absolute character position $pos, optionally
# bolding it, and return the position right after the inserted text.
# (Using a fresh zero-width Range for every insert keeps the math
# simple -- no paragraph-mark / CR surprises to account for.)
# --------------------------------------------------------------------
function Insert-Run($pos, $text, $bold) {
    $seg = $d.Range($pos, $pos)
    $seg.InsertAfter($text)
    $newPos = $pos + $text.Length
    if ($bold) {
        $segBold = $d.Range($pos, $newPos)
        $segBold.Bold = 1
    }
    return $newPos
}

# ------------------------------------------------------------------
# 1) "After the weather-1.0.war  file is generated..." paragraph:
#    rename weather-1.0.war -> weather.war and make that run bold.
# ------------------------------------------------------------------
$p1 = $d.Paragraphs.Item(14)
$r1 = $p1.Range
$found = $r1.Find.Execute("weather-1.0.war")
if ($found) {
    $r1.Text = "weather.war"
    $r1.Bold = 1
}

# ------------------------------------------------------------------
# 2) Insert a brand-new list paragraph right before the
#    "For Tomcat, copy .war file to webapps directory" paragraph,
#    describing the context root that must be set.
# ------------------------------------------------------------------
$pTomcat = $d.Paragraphs.Item(15)
$pTomcat.Range.InsertParagraphBefore()

$pNew = $d.Paragraphs.Item(15)
$pos = $pNew.Range.Start

$pos = Insert-Run $pos "Make sure the " $false
$pos = Insert-Run $pos "context" $true
$pos = Insert-Run $pos " " $false
$pos = Insert-Run $pos "root" $true

$openQuote = [char]0x2018
$midText = " of the application is set to " + $openQuote
$pos = Insert-Run $pos $midText $false

$pos = Insert-Run $pos "/weather" $true

$closeQuote = [string][char]0x2019
$pos = Insert-Run $pos $closeQuote $false

# ------------------------------------------------------------------
# 3) Rewrite the "For Tomcat, copy .war file to webapps directory"
#    paragraph (now shifted down to index 16) so it references the
#    new weather.war name, and restore the _GoBack bookmark in the
#    middle of the sentence.
# ------------------------------------------------------------------
$pTomcat = $d.Paragraphs.Item(16)
$rTomcat = $pTomcat.Range
$clearRange = $d.Range($rTomcat.Start, $rTomcat.End - 1)
$clearRange.Text = ""

$pos = $pTomcat.Range.Start

$pos = Insert-Run $pos "For Tomcat, copy " $false
$pos = Insert-Run $pos "weather" $false
$pos = Insert-Run $pos ".war " $false
$pos = Insert-Run $pos " " $false

$bmRange = $d.Range($pos, $pos)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

$pos = Insert-Run $pos "file to webapps directory" $false

# ------------------------------------------------------------------
# 4) Drop the stale lastRenderedPageBreak hint in front of
#    "Application screen shot" by re-asserting the run text.
# ------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -eq "Application screen shot") {
        $p.Range.Text = "Application screen shot"
        break
    }
}

# ------------------------------------------------------------------
# 5) The old _GoBack bookmark (previously sitting after the final
#    picture) has already been relocated by the Bookmarks.Add call
#    above -- Word only ever keeps a single _GoBack bookmark and
#    simply moves it, so there is nothing further to clean up here.
# ------------------------------------------------------------------
